# Fix typos in two slide titles:
#   Slide 5 (title "Pre-LabWrap-up", split across 3 runs) -> "Pre-Lab Wrap-up"
#   Slide 6 (title "Week 3 Lab: Course Scedule", split across 2 runs) -> "Week 3 Lab: Course Schedule"
#
# Simply assigning TextRange.Text preserves the existing run boundaries
# (word-level diff), which would leave the old multi-run split in place.
# Deleting the text first and re-inserting it collapses it into a single
# run (inheriting the formatting of the original first run), matching
# how the author actually fixed the typo.

$p = $ppt.ActivePresentation

$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(1).TextFrame.TextRange
$tr5.Delete()
$null = $tr5.InsertAfter("Pre-Lab Wrap-up")

$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(1).TextFrame.TextRange
$tr6.Delete()
$null = $tr6.InsertAfter("Week 3 Lab: Course Schedule")
